$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15545
$ws1.Range("F9").Value = 15395
$ws1.Range("F11").Value = 8990
$ws1.Range("F27").Value = 23
$ws1.Range("F32").Value = 59
$ws1.Range("F35").Value = 315
$ws1.Range("F36").Value = 450
$ws1.Range("F38").Value = 5512

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15545
$ws4.Range("F9").Value = 15395
$ws4.Range("F11").Value = 8990
$ws4.Range("F27").Value = 23
$ws4.Range("F34").Value = 59
$ws4.Range("F37").Value = 315
$ws4.Range("F38").Value = 450
$ws4.Range("F40").Value = 5512
